$d = $word.ActiveDocument
$t = $d.Tables(1)
$br = [char]11

$cell = $t.Cell(1, 1)
$cell.Range.Text = "37 x 75" + $br + "  7    5" + $br + "  ----" + $br + "3|    |" + $br + "7|    |"

$cell = $t.Cell(1, 2)
$cell.Range.Text = "65 x 25" + $br + "  2    5" + $br + "  ----" + $br + "6|    |" + $br + "5|    |"

$cell = $t.Cell(1, 3)
$cell.Range.Text = "71 x 93" + $br + "  9    3" + $br + "  ----" + $br + "7|    |" + $br + "1|    |"

$cell = $t.Cell(2, 1)
$cell.Range.Text = "25 x 44" + $br + "  4    4" + $br + "  ----" + $br + "2|    |" + $br + "5|    |"

$cell = $t.Cell(2, 2)
$cell.Range.Text = "74 x 87" + $br + "  8    7" + $br + "  ----" + $br + "7|    |" + $br + "4|    |"

$cell = $t.Cell(2, 3)
$cell.Range.Text = "71 x 44" + $br + "  4    4" + $br + "  ----" + $br + "7|    |" + $br + "1|    |"

$cell = $t.Cell(3, 1)
$cell.Range.Text = "56 x 95" + $br + "  9    5" + $br + "  ----" + $br + "5|    |" + $br + "6|    |"

$cell = $t.Cell(3, 2)
$cell.Range.Text = "33 x 49" + $br + "  4    9" + $br + "  ----" + $br + "3|    |" + $br + "3|    |"

$cell = $t.Cell(3, 3)
$cell.Range.Text = "23 x 87" + $br + "  8    7" + $br + "  ----" + $br + "2|    |" + $br + "3|    |"

$cell = $t.Cell(4, 1)
$cell.Range.Text = "20 x 97" + $br + "  9    7" + $br + "  ----" + $br + "2|    |" + $br + "0|    |"

$cell = $t.Cell(4, 2)
$cell.Range.Text = "54 x 96" + $br + "  9    6" + $br + "  ----" + $br + "5|    |" + $br + "4|    |"

$cell = $t.Cell(4, 3)
$cell.Range.Text = "25 x 14" + $br + "  1    4" + $br + "  ----" + $br + "2|    |" + $br + "5|    |"

$cell = $t.Cell(5, 1)
$cell.Range.Text = "70 x 86" + $br + "  8    6" + $br + "  ----" + $br + "7|    |" + $br + "0|    |"

$cell = $t.Cell(5, 2)
$cell.Range.Text = "69 x 72" + $br + "  7    2" + $br + "  ----" + $br + "6|    |" + $br + "9|    |"

$cell = $t.Cell(5, 3)
$cell.Range.Text = "12 x 51" + $br + "  5    1" + $br + "  ----" + $br + "1|    |" + $br + "2|    |"
